$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "145.51") would be
# auto-converted to a numeric type by COM Value assignment. The source file
# stores these as literal text, so we force text entry (NumberFormat "@"),
# then restore the cell style to Normal so no stray number-format is left
# behind (matches original unstyled cells).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '62.277.83'
$ws.Range("E2").Value = '  -1.99%  '
$ws.Range("D3").Value = '3.005.83'
$ws.Range("E3").Value = '  -2.03%  '
$ws.Range("E4").Value = '  +0.09%  '
Set-TextValue "D5" '582.02'
$ws.Range("E5").Value = '  -1.94%  '
Set-TextValue "D6" '145.51'
$ws.Range("E6").Value = '  -5.69%  '
$ws.Range("E7").Value = '  +0.03%  '
Set-TextValue "D8" '0.528'
$ws.Range("E8").Value = '  -2.45%  '
$ws.Range("D9").Value = '3.006.82'
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("E10").Value = '  -5.29%  '
Set-TextValue "D11" '5.77'
$ws.Range("E11").Value = '  -1.46%  '
Set-TextValue "D12" '0.463'
$ws.Range("E12").Value = '  +2.63%  '
$ws.Range("E13").Value = '  -4.03%  '
Set-TextValue "D14" '34.44'
$ws.Range("E14").Value = '  -6.71%  '
$ws.Range("E15").Value = '  +2.22%  '
$ws.Range("D16").Value = '3.500.40'
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("E17").Value = '  -1.31%  '
$ws.Range("D18").Value = '62.275.93'
$ws.Range("E18").Value = '  -1.92%  '
$ws.Range("D19").Value = '3.007.89'
$ws.Range("E19").Value = '  -1.93%  '
Set-TextValue "D20" '454.06'
$ws.Range("E20").Value = '  -7.46%  '
Set-TextValue "D21" '13.95'
$ws.Range("E21").Value = '  -3.39%  '
Set-TextValue "D22" '0.687'
$ws.Range("E22").Value = '  -2.86%  '
$ws.Range("E23").Value = '  -2.12%  '
Set-TextValue "D24" '81.61'
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("E25").Value = '  -4.14%  '
$ws.Range("E26").Value = '  -10.35%  '
$ws.Range("E27").Value = '  +0.05%  '
Set-TextValue "D28" '9.99'
$ws.Range("E28").Value = '  -7.20%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  -3.23%  '
$ws.Range("E31").Value = '  -6.10%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D32" '2.09'
$ws.Range("E32").Value = '  -5.74%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D33" '28.34'
$ws.Range("E33").Value = '  +3.66%  '
$ws.Range("E34").Value = '  -2.97%  '
$ws.Range("D35").Value = '0.0₃0797'
$ws.Range("E35").Value = '  -3.16%  '
$ws.Range("E36").Value = '  -3.27%  '
$ws.Range("E37").Value = '  -3.95%  '
$ws.Range("E38").Value = '  -5.33%  '
$ws.Range("E39").Value = '  -1.23%  '
Set-TextValue "D40" '50.21'
$ws.Range("E40").Value = '  -0.82%  '
$ws.Range("E41").Value = '  -13.27%  '
Set-TextValue "D42" '0.118'
$ws.Range("E42").Value = '  +3.67%  '
Set-TextValue "D43" '390.15'
$ws.Range("E43").Value = '  -10.94%  '
$ws.Range("E44").Value = '  -1.86%  '
$ws.Range("E45").Value = '  -7.83%  '
$ws.Range("D46").Value = '2.720.24'
$ws.Range("E46").Value = '  -4.28%  '
Set-TextValue "D47" '37.09'
$ws.Range("E47").Value = '  -5.45%  '
Set-TextValue "D48" '128.61'
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue "D50" '2.20'
$ws.Range("E50").Value = '  -2.00%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D51" '0.109'
$ws.Range("E51").Value = '  -0.96%  '
